$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Source:" header (bold, same style as the other section headers
# like A13 "Category Definitions") plus the two source detail lines.
$ws.Range("A18").Value = "Source:"
$ws.Range("A18").Font.Bold = $true
$ws.Range("A19").Value = "National Corrections Reporting Program"
$ws.Range("A20").Value = "https://www.bjs.gov/index.cfm?ty=dcdetail&iid=268"

# Match the author's final selection (cell A18 was left selected).
[void]$ws.Range("A18").Select()
